$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = '[''HMV23/000055/0923/2'', ''HMV23/000055/0923/6'', ''HMV23/000059/0923/5'', ''HMV23/000043/0823/12'']'
$ws.Range("C4").Value = '[''DURING INSPECTION OBSERVED LH WING SLAT#5 WEATHER SEAL ERODED.SAME TO BE REPLACED.'', ''DURING INSPECTION, OBSERVED RH SIDE WING SLAT NO#5 WEATHER SEAL ERODED.SAME TO BE REPLACED.'', ''DURING ARRIVAL INSPECTION , OBSERVED FOLLOWINS DICREPANCIES.\n\n1) RH SIDE WING NO#1 SLAT, NO.3 TRACK SQURE SEAL UPPER BULB SEAL DEGRADED.\n2) RH SIIDE WING SLAT NO#5 OUT BOARD WEATHER SEAL ERODED. SAME TO BE REPLACED.'', ''WHILE ARRIVAL INSPECTION OBSERVED RH WING SLAT#5 OUTBOARD WEATHER SEAL ERODED.SAME TO BE REPLACED.'']'
$ws.Range("D4").Value = '[2.0, 3.0, 0.0, 2.0]'
$ws.Range("B7").Value = '[''HMV23/000059/0923/8'', ''HMV23/000059/0923/9'', ''HMV23/000055/0923/14'', ''HMV23/000055/0923/13'']'
$ws.Range("C7").Value = '[''DURING ARRIVAL INSPECTION OF FORWARD CARGO , FOUND HARDWARE DAMAGED/MISSING'', ''DURING ARRIVAL INSPECTION OF AFT CARGO , FOUND \n1) CARGO SECTION HOSE MISSING/ DAMAGED \n2) FOUND HARDWARE DAMAGED/MISSING'', ''DURING ARRIVAL INSPECTION OF AFT CARGO , FOUND \n1) CARGO SECTION HOSE MISSING/ DAMAGED \n2) FOUND HARDWARE DAMAGED/MISSING'', ''DURING ARRIVAL INSPECTION OF FORWARD CARGO , FOUND \n1) CARGO SECTION HOSE MISSING/ DAMAGED \n2) FOUND HARDWARE DAMAGED/MISSING'']'
$ws.Range("D7").Value = '[5.0, 7.0, 5.0, 5.0]'
$ws.Range("B8").Value = '[''HMV23/000059/0923/10'', ''HMV23/000055/0923/22'']'
$ws.Range("C8").Value = '[''DURING OPPORTUNITY INSPECTION FOUND BELLY PANEL HARDWARE DAMAGED/MISSING'', ''DURING OPP0RTUNITY INSPECTION FOUND BELLY PANEL HARDWARE DAMAGED/MISSING'']'
$ws.Range("D8").Value = '[2.0, 4.0]'
$ws.Range("B9").Value = '[''HMV23/000055/0923/15'', ''HMV23/000059/0923/11'']'
$ws.Range("C9").Value = '[''DURING INSPECTION FOUND PISTON TYPE DRAIN VALVES FOUND DAMAGED/ DEFECTIVE'', ''DURING INSPECTION FOUND PISTON TYPE DRAIN VALVES FOUND DAMAGED/ DEFECTIVE \n\nVALVE P/N ABS0341-2-01'']'
$ws.Range("D9").Value = '[3.0, 4.0]'
$ws.Range("A11").Value = '[''AWR23/003086/0923'', ''AIRCRAFT ARRIVAL INSPECTION'']'
$ws.Range("B13").Value = '[''HMV23/000055/0923/25'', ''HMV23/000059/0923/17'']'
$ws.Range("C13").Value = '[''DURING INSPECTION OF PAX SEATS FOLLOWING FINDINGS OBSERVED:-\n1. BAGGAGE BAR HARDWARES FOUND MISSING AT MULTIPLE LOCATIONS.\n2. SEAT BELTS QTY 17 FOUND FRAYED AND WORN OUT.\n3. ARM CAP QTY 8 FOUND DAMAGED.\n4. SEAT FLANK FAIRINGS AND SPREADER FAIRING FOUND DAMAGED AT MULTIPLE LOCATIONS.'', ''DURING INSPECTION OF PAX SEATS FOLLOWING FINDINGS OBSERVED:-\n1. BAGGAGE BAR HARDWARES FOUND MISSING AT MULTIPLE LOCATIONS.\n2. SEAT BELTS QTY 12 FOUND FRAYED AND WORN OUT.\n3. ARM CAP QTY 1 FOUND DAMAGED.\n4. SEAT FLANK FAIRINGS AND SPREADER FAIRING FOUND DAMAGED AT MULTIPLE LOCATIONS.'']'
$ws.Range("D13").Value = '[8.0, 16.0]'
$ws.Range("B14").Value = '[''HMV23/000055/0923/1'', ''HMV23/000055/0923/9'', ''HMV23/000055/0923/7'', ''HMV23/000055/0923/8'']'
$ws.Range("C14").Value = '[''DURING ARRIVAL INSPECTION , OBSERVED DENT ON SLAT NO #2 LIP AREA AT RH SIDE WING. ASSESSMENT TO BE CARRY OUT.\n\nSLAT NO #2 P/N- D5746092000500\n SA7205698'', ''DURING ARRIVAL INSPECTION, OBSERVED SMALL DENTS ON SLAT NO #1 TOP SKIN AT RH SIDE WING. LOCATION OF DENT 3.5 INCH FROM T/E AND 37 INCH FROM INBD EDGE OF SLAT.\nSLAT NO#1 P/N- D5746091000300\n S/N- SA7205917'', ''DURING ARRIVAL INSPECTION , OBSERVED DENT ON SLAT NO #2 LIP AREA AT RH SIDE WING. DENT LOCATION, 3 INCH FROM T/E AND 67.5 INCH FROM OUTBD EDGE .ASSESSMENT TO BE CARRY OUT.\nSLAT NO #2 P/N- D5746092000500\n S/N- SA7205698'', ''DURING ARRIVAL INSPECTION, OBSERVED SMALL DENTS ON SLAT NO #1 TOP SKIN AT RH SIDE WING. LOCATION OF DENT 3.5 INCH FROM T/E AND 37 INCH FROM INBD EDGE OF SLAT.\nSLAT NO#1 P/N- D5746091000300\n S/N- SA7205917'']'
$ws.Range("D14").Value = '[2.0, 7.0, 7.0, 3.0]'
$ws.Range("B16").Value = '[''HMV23/000055/0923/20'', ''HMV23/000055/0923/19'']'
$ws.Range("C16").Value = '[''1.CARRY OUT HFEC INSPECTION AS PER NTM 51-10-08 ON DENT AND SURROUNDING AREA OF DENT ON RH WING SLAT#1 NOSE TOP SKIN TO ENSURE CRACK FREE CONDITION EXISTS AS REFERRED BY SRM TASK 57-40-00-283-007 REV:37 DTD:AUG 01/2023.\n\nDAMAGE DIMENSION AND LOCATION ARE\n\nDENT ON NOSE TOP SKIN OF SLAT #1 . DAMAGE DIMENSION \nA) LENGTH: 28MM\n\nB) WIDTH: 25MM\n\nC) DEPTH: 0.30MM\n\nD) LOCATION: DISTANCE FROM RIB TRACK 2 RIVET ROW CENTER =109MM, DISTANCE FROM UPEER GIRDER RIVET ROW CENTER =25MM'', ''1.CARRY OUT HFEC INSPECTION AS PER NTM 51-10-08 ON DENT AND SURROUNDING AREA OF DENT ON RH WING SLAT#2 NOSE TOP SKIN TO ENSURE CRACK FREE CONDITION EXISTS AS REFERRED BY SRM TASK 57-40-00-283-007 REV:37 DTD:AUG 01/2023.\n\nDAMAGE DIMENSION AND LOCATION ARE\n\nDENT ON NOSE TOP SKIN OF SLAT #2 . DAMAGE DIMENSION \nA) LENGTH: 84MM\n\nB) WIDTH: 35MM\n\nC) DEPTH: 1.35MM\n\nD) LOCATION: DISTANCE FROM RIB WFX5997 FASTENER ROW CENTER =76MM, DISTANCE FROM UPPER GIRDER RIVET ROW CENTER =65MM'']'
$ws.Range("B17").Value = '[''HMV23/000043/0823/5'', ''HMV23/000043/0823/3'']'
$ws.Range("C17").Value = '[''DURING ARRIVAL INSPECTION OBSERVED DELAMINATION ON AFT CARGO FLOOR PANELS 151SG,151DF,151GF, 51LF, 151PF.\n1)REMOVE AFT CARGO FLOOR PANELS 151SG,151DF,151GF, 51LF, 151PF FOR FABRICATION.\n2) INSTALL AFT CARGO FLOOR PANELS 151SG,151DF,151GF, 51LF, 151PF POST FABRICATION.\n3) SEVERAL HARDWARE FOUND MISSING ON AFT CARGO SIDEWALL AND CEILING PANELS.'', ''DURING ARRIVAL INSPECTION OBSERVED DELAMINATION ON FWD CARGO FLOOR PANELS 131AF,131EF,131FF,131JF.\n1)REMOVE FWD CARGO FLOOR PANEL 131AF,131EF,131FF,131JF FOR FABRICATION.\n2) INSTALL FLOOR PANEL 131AF,131EF,131FF,131JF POST FABRICATION.\n3) SEALS MISSING ON SOME FWD CARGO FLOOR PANELS.\n4) SEVERAL HARDWARE FOUND MISSING ON FWD CARGO SIDEWALL AND CEILING PANELS.'']'
$ws.Range("B19").Value = '[''HMV23/000043/0823/8'', ''HMV23/000043/0823/9'']'
$ws.Range("C19").Value = '[''REMOVAL OF LH WING SLAT #2 (P/N-D5746092000400, S/N-SA7202688)TO BE CARRIED OUT FOR REPAIR.'', ''POST REPAIR INSTALLATION OF LH WING SLAT #2(P/N-D5746092000400, S/N-SA7202688) TO BE CARRIED OUT .'']'
$ws.Range("D19").Value = '[5.0, 10.0]'
